# Update the two-digit ÷ one-digit division answers in the table.
# Each "old" text below occurs exactly once in the document, so a
# simple MatchWholeWord-less, case-sensitive Find/Replace on the whole
# cell text is safe. The one exception is that "36÷3=12, 0" is both a
# target replacement value (for "15÷4=3, 3") and a source value (to be
# replaced by "85÷2=42, 1") elsewhere in the document; we therefore
# perform that particular replacement *before* creating the new
# "36÷3=12, 0" text, so Find.Execute can never match the freshly
# written text.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "40÷4=10, 0" "19÷4=4, 3"
Replace-Text "45÷3=15, 0" "12÷3=4, 0"
Replace-Text "10÷2=5, 0" "84÷8=10, 4"
Replace-Text "78÷8=9, 6" "94÷2=47, 0"
Replace-Text "25÷8=3, 1" "44÷9=4, 8"
Replace-Text "53÷2=26, 1" "86÷4=21, 2"
Replace-Text "46÷9=5, 1" "38÷7=5, 3"
Replace-Text "36÷4=9, 0" "92÷2=46, 0"
Replace-Text "68÷3=22, 2" "14÷4=3, 2"
Replace-Text "19÷6=3, 1" "21÷6=3, 3"
Replace-Text "43÷9=4, 7" "78÷9=8, 6"
Replace-Text "26÷5=5, 1" "20÷5=4, 0"
Replace-Text "79÷2=39, 1" "86÷3=28, 2"
Replace-Text "14÷7=2, 0" "91÷9=10, 1"
Replace-Text "57÷8=7, 1" "83÷8=10, 3"
Replace-Text "92÷4=23, 0" "20÷3=6, 2"
Replace-Text "83÷9=9, 2" "17÷5=3, 2"
Replace-Text "46÷2=23, 0" "55÷4=13, 3"
Replace-Text "18÷6=3, 0" "79÷3=26, 1"

# Do this one before re-introducing "36÷3=12, 0" below, since the old
# text here is identical to the new text produced by the next step.
Replace-Text "36÷3=12, 0" "85÷2=42, 1"

Replace-Text "15÷4=3, 3" "36÷3=12, 0"

Replace-Text "55÷6=9, 1" "97÷2=48, 1"
Replace-Text "98÷9=10, 8" "32÷5=6, 2"
Replace-Text "26÷9=2, 8" "18÷9=2, 0"
Replace-Text "90÷4=22, 2" "19÷5=3, 4"
